# Fixed update to excel issue
# 1. Rename "Requested quantity" headers to reflect new naming convention
# 2. Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet (Monthly Trend)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the header formatting (bold, centered, bordered) from the Weekly Quantity header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Reuse the date-column formatting (custom date/time number format) from column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A33").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$forecastData = @(
    @(44983.99999999999, 85, -20.08482791750849, 187.8902386863956),
    @(45004.99999999999, 85, -20.419606311402, 189.0297450551411),
    @(45011.99999999999, 85, -17.50545163153635, 188.8695604787052),
    @(45018.99999999999, 85, -21.96369454344176, 189.0090619142208),
    @(45151.99999999999, 86, -25.68345222352805, 186.7137659941805),
    @(45158.99999999999, 86, -25.66704304867037, 200.1059845697767),
    @(45165.99999999999, 86, -24.85208029716501, 200.5413392038703),
    @(45172.99999999999, 86, -22.18619460046977, 186.7332584201667),
    @(45179.99999999999, 86, -20.30125398234767, 195.1696935940191),
    @(45186.99999999999, 86, -25.48604813612787, 195.9742363044242),
    @(45193.99999999999, 86, -24.53559317414004, 195.9927615296611),
    @(45207.99999999999, 86, -22.80417809597028, 190.2138125272641),
    @(45214.99999999999, 86, -13.61010705188423, 194.8223520717324),
    @(45228.99999999999, 86, -24.03087394872215, 189.7919448399541),
    @(45235.99999999999, 86, -21.24030485546874, 198.4168027059821),
    @(45242.99999999999, 87, -14.78277549568434, 190.5335680106415),
    @(45249.99999999999, 87, -24.12977069875354, 193.9710785831379),
    @(45256.99999999999, 87, -25.07845627308344, 195.2860388470222),
    @(45270.99999999999, 87, -23.30976571605917, 198.2774457079426),
    @(45277.99999999999, 87, -25.88000039605541, 200.0527098499818),
    @(45298.99999999999, 87, -19.18067255507814, 187.0152807689812),
    @(45305.99999999999, 87, -22.3956795381919, 193.5000866508918),
    @(45312.99999999999, 87, -13.02013616494527, 196.0032460771824),
    @(45361.99999999999, 87, -9.604419115750575, 211.3889727134011),
    @(45368.99999999999, 87, -27.0870607558845, 202.2635242358425),
    @(45375.99999999999, 87, -29.73366543700056, 198.9620337950973),
    @(45382.99999999999, 88, -17.98853558095734, 201.3519410382065),
    @(45389.99999999999, 88, -17.43384749360223, 197.2122856971083),
    @(45396.99999999999, 88, -22.27958711740193, 196.1045300556407),
    @(45403.99999999999, 88, -18.6498980600304, 192.9475405444393),
    @(45410.99999999999, 88, -20.48694836832432, 200.0139638155124),
    @(45417.99999999999, 88, -17.05896779955965, 198.7061026667855)
)

$row = 2
foreach ($entry in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $entry[0]
    $wsForecast.Cells.Item($row, 2).Value = $entry[1]
    $wsForecast.Cells.Item($row, 3).Value = $entry[2]
    $wsForecast.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

[void]$wsForecast.Range("A1").Select()
